# Registro: se estan haciendo las validaciones del registro.
# Agrega dos nuevas filas de datos (registros) a la hoja "Hoja1",
# completa la celda G2 que faltaba y deja seleccionada la celda A4,
# tal como queda el archivo tras capturar los nuevos registros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Completar el valor que faltaba en la fila 2 ---
$ws.Range("G2").Value = "No tiene"

# --- Fila 3: nuevo registro "SS" / "SS" / 11 / "SS" / 111 / "SS" / "No tiene" ---
$ws.Range("A3").Value = "SS"
$ws.Range("B3").Value = "SS"
$ws.Range("C3").Value = 11
$ws.Range("D3").Value = "SS"
$ws.Range("E3").Value = 111
$ws.Range("F3").Value = "SS"
$ws.Range("G3").Value = "No tiene"

# --- Fila 4: nuevo registro "dd" / "dd" / 22 / "ds" / "23" / "xb" / "No tiene" ---
$ws.Range("A4").Value = "dd"
$ws.Range("B4").Value = "dd"
$ws.Range("C4").Value = 22
$ws.Range("D4").Value = "ds"

# La columna E de la fila 4 contiene el texto "23" (no el numero 23).
# Se construye mediante una formula de texto y se pega como valor para
# conservarlo como cadena de texto en vez de convertirlo a numero.
$scratch = $ws.Range("Z100")
$scratch.Formula = "=TEXT(23,""0"")"
$scratch.Copy()
$ws.Range("E4").PasteSpecial(-4163)
$scratch.ClearContents()

$ws.Range("F4").Value = "xb"
$ws.Range("G4").Value = "No tiene"

# Selecciona la celda A4, que quedo como celda activa tras capturar el registro.
$ws.Range("A4").Select()
